# "added # of adj. gene airs" -- replace the placeholder histogram bin
# counts in column B (Sheet1!B1:B41) with the real "# of adjacent gene
# pairs" counts. Column A (the PCC bin edges) and the B42 SUM(B1:B41)
# total are left alone; B42 recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New bin counts for B1:B41 (idx 0..40 in the chart's numCache maps to
# rows 1..41 here).
$newCounts = @(
    0,
    824,
    2176,
    1413,
    576,
    243,
    167,
    121,
    208,
    450,
    218,
    199,
    94,
    118,
    271,
    151,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    58,
    138,
    100,
    255,
    78,
    360,
    0,
    227,
    298,
    412,
    586,
    1146,
    1399,
    1605
)

for ($i = 0; $i -lt $newCounts.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = $newCounts[$i]
}
